$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.264.37"
$ws.Range("E2").Value = "  -5.90%  "
$ws.Range("D3").Value = "1.675.77"
$ws.Range("E3").Value = "  -3.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.96"
$ws.Range("E5").Value = "  -3.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5094"
$ws.Range("E6").Value = "  -11.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2663"
$ws.Range("E8").Value = "  -2.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06358"
$ws.Range("E9").Value = "  -3.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.64"
$ws.Range("E10").Value = "  -6.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07374"
$ws.Range("E11").Value = "  -2.21%  "
$ws.Range("D12").Value = "1.679.66"
$ws.Range("E12").Value = "  -3.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.548"
$ws.Range("E13").Value = "  -3.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5811"
$ws.Range("E14").Value = "  -3.47%  "
$ws.Range("D15").Value = "1.899.34"
$ws.Range("E15").Value = "  -3.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008566"
$ws.Range("E16").Value = "  -2.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.19"
$ws.Range("E17").Value = "  -12.55%  "
$ws.Range("D18").Value = "26.353.25"
$ws.Range("E18").Value = "  -5.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.950"
$ws.Range("E19").Value = "  -6.93%  "
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.85"
$ws.Range("E21").Value = "  -3.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "190.21"
$ws.Range("E22").Value = "  -7.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.214"
$ws.Range("E23").Value = "  -5.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.006"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.09"
$ws.Range("E25").Value = "  -4.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.690"
$ws.Range("E26").Value = "  -5.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1174"
$ws.Range("E27").Value = "  -4.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.73"
$ws.Range("E28").Value = "  -2.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05819"
$ws.Range("E29").Value = "  -5.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.280"
$ws.Range("E30").Value = "  -7.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.321"
$ws.Range("E31").Value = "  -4.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.534"
$ws.Range("E32").Value = "  -5.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.519"
$ws.Range("E33").Value = "  -5.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.650"
$ws.Range("E34").Value = "  -1.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.010"
$ws.Range("E35").Value = "  -2.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6003"
$ws.Range("E36").Value = "  -6.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.359"
$ws.Range("E37").Value = "  -2.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.650"
$ws.Range("E38").Value = "  -2.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01617"
$ws.Range("E39").Value = "  -3.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.039"
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("D41").Value = "1.080.59"
$ws.Range("E41").Value = "  -4.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8612"
$ws.Range("E42").Value = "  -1.28%  "
$ws.Range("E43").Value = "  +0.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.01"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").Value = "1.827.08"
$ws.Range("E45").Value = "  -3.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000112"
$ws.Range("E46").Value = "  +3.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.96"
$ws.Range("E47").Value = "  -5.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.007"
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.095"
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4297"
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05192"
$ws.Range("E51").Value = "  -3.47%  "
